$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy current row 8 (GMHO:0000191 / repeated measures study design) down to the
#    new row 9, preserving its values and formatting (style index 2).
$ws.Range("A8:V8").Copy($ws.Range("A9:V9"))

# 2. Overwrite row 8 with the new "planned process" entity (OBI:0000011), keeping
#    the existing s="2" formatting already present on row 8.
$ws.Range("A8").Value = "OBI:0000011"
$ws.Range("B8").Value = "planned process"
$ws.Range("C8").Value = "A process that realizes a plan which is the concretization of a plan specification."
$ws.Range("D8").Value = "process"
$ws.Range("E8").Value = "COB:0000082"

# 3. Row 5 loses its special formatting (style index 2 cleared back to default),
#    matching the look of rows 3/4, and its Curation status changes.
$ws.Range("A5:V5").ClearFormats()
$ws.Range("S5").Value = "External"
